# BOM update:
#  1. Insert a new "Battery" BOM line right after the header row.
#  2. Extend the 4K99/0603 resistor designator list (+R136, R137) and bump its qty.
#  3. Replace the RTC evaluation-board line with the DS3231SN# RTC IC.
#  4. Replace the 74LVX3245 line with the SN74LVC8T245DWR level translator.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- 1. Insert new row 2 (Battery) -----------------------------------
$ws.Rows(2).Insert()

$ws.Range("A2").Value = "Battery"
$ws.Range("B2").Value = "Battery Holder (Open) Coin, 12.0mm 1 Cell SMD (SMT) Tab"
$ws.Range("C2").Value = "BT1"
$ws.Range("D2").Value = "BC501SM-TR-ND"
$ws.Range("E2").Value = "MPD"
$ws.Range("F2").Value = "BC501SM-TR"
$ws.Range("G2").Value = 1

# Re-apply the standard data-row formatting (Insert() leaves the new row
# with a default style, and writing .Value to text cells drops the
# quotePrefix formatting used throughout the table) by pasting formats
# from the row below, which still carries the original look.
$ws.Range("A3:G3").Copy()
$ws.Range("A2:G2").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

# --- 2. 4K99/0603 resistor row (now row 12): extend designators + qty --
$ws.Range("C12").Value = "R7, R9, R10, R11, R125, R126, R127, R128, R129, R130, R131, R132, R133, R136, R137"
$ws.Range("G12").Value = 15

$ws.Range("A12").Copy()
$ws.Range("C12").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

# --- 3. RTC row (now row 18): RTC3013 eval board -> DS3231SN# RTC IC ---
$ws.Range("A18").Value = "DS3231SN#"
$ws.Range("B18").Value = "Real Time Clock (RTC) IC Clock/Calendar I" + [char]0x00B2 + "C, 2-Wire Serial 16-SOIC (0.295"", 7.50mm Width)"
$ws.Range("D18").Value = "DS3231SN#-ND"
$ws.Range("E18").Value = "Maxim"
$ws.Range("F18").Value = "DS3231SN#"

$ws.Range("C18").Copy()
$ws.Range("A18").PasteSpecial($xlPasteFormats)
$ws.Range("B18").PasteSpecial($xlPasteFormats)
$ws.Range("D18").PasteSpecial($xlPasteFormats)
$ws.Range("E18").PasteSpecial($xlPasteFormats)
$ws.Range("F18").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

# --- 4. Bus buffer row (now row 19): 74LVX3245 -> SN74LVC8T245DWR ------
$ws.Range("A19").Value = "SN74LVC8T245DWR"
$ws.Range("B19").Value = "Voltage Level Translator Bidirectional 1 Circuit 8 Channel 24-SOIC"
$ws.Range("D19").Value = "296-23280-1-ND"
$ws.Range("E19").Value = "TI"
$ws.Range("F19").Value = "SN74LVC8T245DWR"

$ws.Range("C19").Copy()
$ws.Range("A19").PasteSpecial($xlPasteFormats)
$ws.Range("B19").PasteSpecial($xlPasteFormats)
$ws.Range("D19").PasteSpecial($xlPasteFormats)
$ws.Range("E19").PasteSpecial($xlPasteFormats)
$ws.Range("F19").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

Write-Output "BOM updated"
